# db_tempi_base.xlsx -- refresh simulated "tempi" (time/cost) figures on Foglio1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- Row 2 : Pikes Peak ---
$ws.Range("B2").Value = 199.17
$ws.Range("C2").Value = 258.72
$ws.Range("D2").Value = 53.8
$ws.Range("E2").Value = 12.29
$ws.Range("F2").Value = 71.11

# --- Row 3 : Rally ---
$ws.Range("B3").Value = 199.17
$ws.Range("C3").Value = 258.72
$ws.Range("D3").Value = 53.8
$ws.Range("E3").Value = 12.29
$ws.Range("F3").Value = 86.54

# --- Row 4 : RS ---
$ws.Range("B4").Value = 196.49
$ws.Range("C4").Value = 258.72
$ws.Range("D4").Value = 53.8
$ws.Range("E4").Value = 12.29
$ws.Range("F4").Value = 73.79

# --- Row 5 : S ---
$ws.Range("B5").Value = 197.49
$ws.Range("C5").Value = 258.72
$ws.Range("D5").Value = 53.8
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 73.79

# --- Row 6 : S Grand Tour ---
$ws.Range("B6").Value = 196.49
$ws.Range("C6").Value = 258.72
$ws.Range("D6").Value = 53.8
$ws.Range("E6").Value = 12.29
$ws.Range("F6").Value = 73.79

# --- Row 7 : Standard ---
$ws.Range("B7").Value = 186.7
$ws.Range("C7").Value = 258.72
$ws.Range("D7").Value = 53.8
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 73.79

# B7 picks up its own number format (2 decimal places) instead of the
# sheet default -- new style entry shows up in cellXfs.
$ws.Range("B7").NumberFormat = "0.00"

# The author's last selection before saving moved to F18.
$ws.Range("F18").Select()

# Reposition the workbook window to match the author's last saved layout.
$excel.ActiveWindow.Left = 5520
$excel.ActiveWindow.Top = 3240
